# Update cryptos list data (Price and Volume(1h) columns) as scraped on
# Mon Mar 25 02:45:55 UTC 2024 by GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new numeric-looking price would otherwise be auto-parsed
# as a number by Excel are first forced to Text format so the literal
# string (e.g. trailing zero in "573.80") survives exactly as scraped.

$ws.Range("D2").Value = "66.677.09"
$ws.Range("E2").Value = "  +3.59%  "
$ws.Range("D3").Value = "3.439.76"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.80"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.27"
$ws.Range("E6").Value = "  +5.69%  "
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").Value = "3.431.76"
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.644"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.97"
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.40"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "3.987.15"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "3.442.68"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "66.701.10"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("E21").Value = "  +2.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.82"
$ws.Range("E22").Value = "  +8.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.42"
$ws.Range("E23").Value = "  +16.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.03"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.26"
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.71"
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("E29").Value = "  +5.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.40"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  +7.79%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "593.95"
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.72"
$ws.Range("E34").Value = "  +4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.26"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").Value = "3.175.40"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("E45").Value = "  +4.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.79"
$ws.Range("E46").Value = "  +21.63%  "
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.20"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.76"
$ws.Range("E49").Value = "  +7.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "140.36"
$ws.Range("E51").Value = "  -1.07%  "
